$d = $word.ActiveDocument

$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host ("WARNING: replace failed for: " + $old)
    }
}

# --- Simple numeric updates in body paragraphs ---
Replace-Text "Number of tokens: 16130" "Number of tokens: 15477"
Replace-Text "Average number of words per sentence: 18.44" "Average number of words per sentence: 18.89"

# --- Token bigrams / trigrams (filtered stray \n artifact tokens) ---
Replace-Text ": 237, \n The : 98, of the : 82" ": 240, of the : 82, , and : 67"
Replace-Text ": 40, \\ "" \n : 35, \\ "" . : 28" ": 40, \\ "" . : 28, \\ "" , : 13"

# --- Fine-grained / Universal POS n-gram counts ---
Replace-Text "DT NN: 671, NNP NNP: 611, IN DT: 587" "DT NN: 671, NNP NNP: 608, IN DT: 586"
Replace-Text "IN DT NN: 293, NNP NNP NNP: 201, DT NN IN: 195" "IN DT NN: 292, NNP NNP NNP: 200, DT NN IN: 195"
Replace-Text "DET NOUN: 780, NOUN PUNCT: 769, NOUN ADP: 706" "DET NOUN: 780, NOUN PUNCT: 780, NOUN ADP: 710"
Replace-Text "ADP DET NOUN: 302, NOUN ADP DET: 242, VERB DET NOUN: 228" "ADP DET NOUN: 302, NOUN ADP DET: 243, VERB DET NOUN: 228"

# --- Number of named entities ---
Replace-Text "Number of named entities: 1627" "Number of named entities: 1614"

# --- POS frequency table (Table 1) ---
$t = $d.Tables.Item(1)

# Row 3 (NNP): occurrences 2060 -> 2017; tokens column \\, US, President -> , US, President
$t.Cell(3, 3).Range.Text = "2017"
$t.Cell(3, 5).Range.Text = ", US, President"

# Row 6 (JJ): occurrences 868 -> 869
$t.Cell(6, 3).Range.Text = "869"

# Row 7 (NNS): occurrences 774 -> 779
$t.Cell(7, 3).Range.Text = "779"

# Row 8 (,  Punct): relative frequency 0.04 -> 0.05
$t.Cell(8, 4).Range.Text = "0.05"

# Row 11 (was VBN -> becomes VBD, new values)
$t.Cell(11, 1).Range.Text = "VBD"
$t.Cell(11, 3).Range.Text = "451"
$t.Cell(11, 5).Range.Text = "said, reported, told"
$t.Cell(11, 6).Range.Text = "plated"

# Row 10 (was _SP/Space/653/... -> becomes VBN/Verb/454/0.03/accused known killed/hospitalised)
$t.Cell(10, 1).Range.Text = "VBN"
$t.Cell(10, 2).Range.Text = "Verb"
$t.Cell(10, 3).Range.Text = "454"
$t.Cell(10, 4).Range.Text = "0.03"
$t.Cell(10, 5).Range.Text = "accused, known, killed"
$t.Cell(10, 6).Range.Text = "hospitalised"

# --- Mark the screenshot image run as NoProofing (adds <w:noProof/>) ---
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = 1

Write-Host "done"
